$wb = $excel.ActiveWorkbook

# Update "Latest Handback DateTime" (column K) for row 2 on the zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-09-09 07:12:29"

# Update "Latest Handback DateTime" (column K) for row 2 on the de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-09 07:12:47"
